# Generate Report for Handoff
# Inserts a new row for source file "29f631d0-13d9-4325-bd0a-8b9626ddc855.md"
# (status "Ready for handoff") between the existing "303ebb0e-...md" row and
# the "41f722b2-...md" row on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A:G, insert new row at worksheet row 3
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$wsOv.Rows.Item(3).Insert()
$loOv.Resize($wsOv.Range("A1:G4"))

$wsOv.Range("A3").Value = "29f631d0-13d9-4325-bd0a-8b9626ddc855.md"
$wsOv.Range("B3").Value = "e2e\29f631d0-13d9-4325-bd0a-8b9626ddc855.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = ""
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-11-09 06:49:49"

# Rebuild hyperlinks for this sheet (row-insert does not shift them)
$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/303ebb0e-a951-4d20-b4a7-c3c5058f5695.md", "", "", "e2e\303ebb0e-a951-4d20-b4a7-c3c5058f5695.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/29f631d0-13d9-4325-bd0a-8b9626ddc855.md", "", "", "e2e\29f631d0-13d9-4325-bd0a-8b9626ddc855.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md", "", "", "e2e\41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": columns A:P, insert new row at worksheet row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$wsZh.Rows.Item(3).Insert()
$loZh.Resize($wsZh.Range("A1:P4"))

$wsZh.Range("A3").Value = "29f631d0-13d9-4325-bd0a-8b9626ddc855.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "29f631d0-13d9-4325-bd0a-8b9626ddc855.d96602369b5deaf6d97f6c87ebe3abe6a746005c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-11-09 06:49:35"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# Rebuild hyperlinks for this sheet (row-insert does not shift them)
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/303ebb0e-a951-4d20-b4a7-c3c5058f5695.md", "", "", "303ebb0e-a951-4d20-b4a7-c3c5058f5695.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e2880f13dd5dfe07cd125ca6084a52b8a7b9f318/e2e/303ebb0e-a951-4d20-b4a7-c3c5058f5695.md", "", "", "303ebb0e-a951-4d20-b4a7-c3c5058f5695.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/29f631d0-13d9-4325-bd0a-8b9626ddc855.md", "", "", "29f631d0-13d9-4325-bd0a-8b9626ddc855.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md", "", "", "41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3b393a719a17017a7c97abe2f0ec646876acb600/e2e/41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md", "", "", "41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md")

# ---------------------------------------------------------------------------
# Sheet "de-de": columns A:P, insert new row at worksheet row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$wsDe.Rows.Item(3).Insert()
$loDe.Resize($wsDe.Range("A1:P4"))

$wsDe.Range("A3").Value = "29f631d0-13d9-4325-bd0a-8b9626ddc855.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "29f631d0-13d9-4325-bd0a-8b9626ddc855.d96602369b5deaf6d97f6c87ebe3abe6a746005c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-11-09 06:49:49"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

# Rebuild hyperlinks for this sheet (row-insert does not shift them)
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/303ebb0e-a951-4d20-b4a7-c3c5058f5695.md", "", "", "303ebb0e-a951-4d20-b4a7-c3c5058f5695.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/375e7d6a290b0f98ffe497a6e71088b2754be419/e2e/303ebb0e-a951-4d20-b4a7-c3c5058f5695.md", "", "", "303ebb0e-a951-4d20-b4a7-c3c5058f5695.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/29f631d0-13d9-4325-bd0a-8b9626ddc855.md", "", "", "29f631d0-13d9-4325-bd0a-8b9626ddc855.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/07a222f3a4bb8f1ca3ecad1e0261961dde199f13/e2e/41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md", "", "", "41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b0e836bc99eeb4201e670544df9fcf63db32966c/e2e/41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md", "", "", "41f722b2-4d1c-4098-9b73-ec2f6844ab2b.md")
